# Auto-generated edit script: refresh Universalis market-price snapshot values
# for the Tonberry_Profits leve-crafting profit sheets (H:N columns).
# Values below are the latest scheduled-runner market data pull;
# columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 877.4231
$ws.Range("I80").Value = 815.8125
$ws.Range("J80").Value = 976
$ws.Range("K80").Value = 2447.4375
$ws.Range("L80").Value = 2928
$ws.Range("M80").Value = -1449.4375
$ws.Range("N80").Value = -4924
# Row 83
$ws.Range("H83").Value = 877.4231
$ws.Range("I83").Value = 815.8125
$ws.Range("J83").Value = 976
$ws.Range("K83").Value = 7342.3125
$ws.Range("L83").Value = 8784
$ws.Range("M83").Value = -2350.3125
$ws.Range("N83").Value = -18768
# Row 86
$ws.Range("H86").Value = 2498.5
$ws.Range("I86").Value = 998.2
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 998.2
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = 124.8
$ws.Range("N86").Value = -12246
# Row 89
$ws.Range("H89").Value = 2498.5
$ws.Range("I89").Value = 998.2
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 4991
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = 625
$ws.Range("N89").Value = -61232
# Row 106
$ws.Range("H106").Value = 2703.3333
$ws.Range("I106").Value = 1342.8
$ws.Range("K106").Value = 1342.8
$ws.Range("M106").Value = -711.8
# Row 129
$ws.Range("H129").Value = 893.5
$ws.Range("J129").Value = 882.2895
$ws.Range("L129").Value = 2646.8685
$ws.Range("N129").Value = -12646.8685
# Row 132
$ws.Range("H132").Value = 1194.1072
$ws.Range("I132").Value = 992.96
$ws.Range("K132").Value = 2978.88
$ws.Range("M132").Value = -448.8800000000001
# Row 138
$ws.Range("H138").Value = 2358.9285
$ws.Range("I138").Value = 2237.6
$ws.Range("J138").Value = 2662.25
$ws.Range("K138").Value = 6712.799999999999
$ws.Range("L138").Value = 7986.75
$ws.Range("M138").Value = -1572.799999999999
$ws.Range("N138").Value = -18266.75
# Row 141
$ws.Range("H141").Value = 3621.625
$ws.Range("I141").Value = 2947.1177
$ws.Range("J141").Value = 5259.7144
$ws.Range("K141").Value = 8841.3531
$ws.Range("L141").Value = 15779.1432
$ws.Range("M141").Value = -3661.3531
$ws.Range("N141").Value = -26139.1432

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3801.5305
$ws.Range("I32").Value = 2157.7073
$ws.Range("J32").Value = 12226.125
$ws.Range("K32").Value = 2157.7073
$ws.Range("L32").Value = 12226.125
$ws.Range("M32").Value = -1870.7073
$ws.Range("N32").Value = -12800.125
# Row 61
$ws.Range("H61").Value = 3281.3225
$ws.Range("I61").Value = 2194.1667
$ws.Range("K61").Value = 2194.1667
$ws.Range("M61").Value = -1982.1667
# Row 132
$ws.Range("H132").Value = 1715.5227
$ws.Range("I132").Value = 1206.4828
$ws.Range("K132").Value = 3619.4484
$ws.Range("M132").Value = -1089.4484
# Row 136
$ws.Range("H136").Value = 3281.3225
$ws.Range("I136").Value = 2194.1667
$ws.Range("K136").Value = 6582.500100000001
$ws.Range("M136").Value = -4032.500100000001

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 254732.25
$ws.Range("I86").Value = 9283.666999999999
$ws.Range("J86").Value = 402001.4
$ws.Range("K86").Value = 9283.666999999999
$ws.Range("L86").Value = 402001.4
$ws.Range("M86").Value = -8160.666999999999
$ws.Range("N86").Value = -404247.4
# Row 89
$ws.Range("H89").Value = 254732.25
$ws.Range("I89").Value = 9283.666999999999
$ws.Range("J89").Value = 402001.4
$ws.Range("K89").Value = 46418.335
$ws.Range("L89").Value = 2010007
$ws.Range("M89").Value = -40802.335
$ws.Range("N89").Value = -2021239
# Row 134
$ws.Range("H134").Value = 5036.206
$ws.Range("I134").Value = 5441.033
$ws.Range("K134").Value = 16323.099
$ws.Range("M134").Value = -13788.099

$ws = $wb.Worksheets.Item("CRP")
# Row 19
$ws.Range("H19").Value = 953.3333
$ws.Range("I19").Value = 927.5
$ws.Range("K19").Value = 927.5
$ws.Range("M19").Value = -757.5
# Row 24
$ws.Range("H24").Value = 953.3333
$ws.Range("I24").Value = 927.5
$ws.Range("K24").Value = 927.5
$ws.Range("M24").Value = -757.5
# Row 31
$ws.Range("H31").Value = 2260.5186
$ws.Range("I31").Value = 1821.9166
$ws.Range("J31").Value = 2611.4
$ws.Range("K31").Value = 1821.9166
$ws.Range("L31").Value = 2611.4
$ws.Range("M31").Value = -1526.9166
$ws.Range("N31").Value = -3201.4
# Row 34
$ws.Range("H34").Value = 2260.5186
$ws.Range("I34").Value = 1821.9166
$ws.Range("J34").Value = 2611.4
$ws.Range("K34").Value = 1821.9166
$ws.Range("L34").Value = 2611.4
$ws.Range("M34").Value = -1619.9166
$ws.Range("N34").Value = -3015.4
# Row 99
$ws.Range("H99").Value = 2656.889
$ws.Range("I99").Value = 1999.5
$ws.Range("K99").Value = 1999.5
$ws.Range("M99").Value = -501.5
# Row 107
$ws.Range("H107").Value = 537.94446
$ws.Range("I107").Value = 493.53333
$ws.Range("K107").Value = 493.53333
$ws.Range("M107").Value = 1426.46667
# Row 126
$ws.Range("H126").Value = 2656.889
$ws.Range("I126").Value = 1999.5
$ws.Range("K126").Value = 5998.5
$ws.Range("M126").Value = -3528.5
# Row 132
$ws.Range("H132").Value = 1955.4048
$ws.Range("I132").Value = 1248.1923
$ws.Range("J132").Value = 3104.625
$ws.Range("K132").Value = 3744.5769
$ws.Range("L132").Value = 9313.875
$ws.Range("M132").Value = -1214.5769
$ws.Range("N132").Value = -14373.875

$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 1033.9333
$ws.Range("J122").Value = 1132.7273
$ws.Range("L122").Value = 10194.5457
$ws.Range("N122").Value = -15094.5457
# Row 131
$ws.Range("H131").Value = 748.34
$ws.Range("J131").Value = 777.1429000000001
$ws.Range("L131").Value = 2331.4287
$ws.Range("N131").Value = -12411.4287
# Row 140
$ws.Range("H140").Value = 1888.0322
$ws.Range("I140").Value = 1101.0625
$ws.Range("K140").Value = 3303.1875
$ws.Range("M140").Value = 1876.8125

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 67.588234
$ws.Range("I2").Value = 15.8
$ws.Range("J2").Value = 89.166664
$ws.Range("K2").Value = 15.8
$ws.Range("L2").Value = 89.166664
$ws.Range("M2").Value = 97.2
$ws.Range("N2").Value = -315.166664
# Row 107
$ws.Range("H107").Value = 1380.6
$ws.Range("I107").Value = 150
$ws.Range("J107").Value = 2201
$ws.Range("K107").Value = 150
$ws.Range("L107").Value = 2201
$ws.Range("M107").Value = 1770
$ws.Range("N107").Value = -6041
# Row 122
$ws.Range("H122").Value = 1590.7826
$ws.Range("I122").Value = 1298.9412
$ws.Range("J122").Value = 2417.6667
$ws.Range("K122").Value = 3896.8236
$ws.Range("L122").Value = 7253.000100000001
$ws.Range("M122").Value = -1446.8236
$ws.Range("N122").Value = -12153.0001
# Row 132
$ws.Range("H132").Value = 3500512
$ws.Range("I132").Value = 5497392.5
$ws.Range("K132").Value = 16492177.5
$ws.Range("M132").Value = -16489647.5

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 3698
$ws.Range("I16").Value = 5582.8335
$ws.Range("J16").Value = 2441.4443
$ws.Range("K16").Value = 5582.8335
$ws.Range("L16").Value = 2441.4443
$ws.Range("M16").Value = -5412.8335
$ws.Range("N16").Value = -2781.4443
# Row 22
$ws.Range("H22").Value = 3037.4443
$ws.Range("J22").Value = 1997.4
$ws.Range("L22").Value = 1997.4
$ws.Range("N22").Value = -2587.4
# Row 27
$ws.Range("H27").Value = 3037.4443
$ws.Range("J27").Value = 1997.4
$ws.Range("L27").Value = 1997.4
$ws.Range("N27").Value = -2211.4
# Row 40
$ws.Range("H40").Value = 2257.1428
$ws.Range("I40").Value = 2300
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 2300
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -2164
$ws.Range("N40").Value = -2272
# Row 122
$ws.Range("H122").Value = 17705
$ws.Range("I122").Value = 21002
$ws.Range("K122").Value = 63006
$ws.Range("M122").Value = -60556
# Row 132
$ws.Range("H132").Value = 2974.72
$ws.Range("I132").Value = 2581
$ws.Range("J132").Value = 3284.0715
$ws.Range("K132").Value = 7743
$ws.Range("L132").Value = 9852.2145
$ws.Range("M132").Value = -5213
$ws.Range("N132").Value = -14912.2145
# Row 136
$ws.Range("H136").Value = 3044
$ws.Range("I136").Value = 2341.8125
$ws.Range("J136").Value = 3793
$ws.Range("K136").Value = 7025.4375
$ws.Range("L136").Value = 11379
$ws.Range("M136").Value = -4475.4375
$ws.Range("N136").Value = -16479

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = ""
$ws.Range("N54").Value = ""
# Row 62
$ws.Range("H62").Value = 2500.6667
$ws.Range("J62").Value = 5003
$ws.Range("L62").Value = 5003
$ws.Range("N62").Value = -6251
# Row 65
$ws.Range("H65").Value = 2500.6667
$ws.Range("J65").Value = 5003
$ws.Range("L65").Value = 25015
$ws.Range("N65").Value = -31255
# Row 107
$ws.Range("H107").Value = 1106.625
$ws.Range("J107").Value = 1142
$ws.Range("L107").Value = 3426
$ws.Range("N107").Value = -7266
# Row 122
$ws.Range("H122").Value = 98447.625
$ws.Range("I122").Value = 112297.29
$ws.Range("K122").Value = 336891.87
$ws.Range("M122").Value = -334441.87
# Row 132
$ws.Range("H132").Value = 1675.4231
$ws.Range("I132").Value = 1142.9445
$ws.Range("J132").Value = 2873.5
$ws.Range("K132").Value = 3428.8335
$ws.Range("L132").Value = 8620.5
$ws.Range("M132").Value = -898.8335000000002
$ws.Range("N132").Value = -13680.5
# Row 136
$ws.Range("H136").Value = 15016879
$ws.Range("I136").Value = 20577616
$ws.Range("J136").Value = 2887.7
$ws.Range("K136").Value = 61732848
$ws.Range("L136").Value = 8663.099999999999
$ws.Range("M136").Value = -61730298
$ws.Range("N136").Value = -13763.1

